# Updates the cryptos list data in the active worksheet.
# Applies per-cell value updates for Price (D) and Volume(1h) (E) columns,
# plus the reordering of a few coin rows (Toncoin/InjectiveProtocol,
# ApeXProtocol/NEARProtocol) and replacement of SEI with TrustWalletToken.
# Cells whose new Price text looks like a plain number are forced to
# Text format first so Excel keeps the exact string (matching the
# source data's literal formatting, e.g. trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.289.59"
$ws.Range("E2").Value = "  -1.30%  "
# Row 3
$ws.Range("D3").Value = "2.770.15"
$ws.Range("E3").Value = "  -0.78%  "
# Row 4
$ws.Range("E4").Value = "  +0.09%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.60"
$ws.Range("E5").Value = "  -3.30%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.73"
$ws.Range("E6").Value = "  -2.02%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.548"
$ws.Range("E7").Value = "  -2.21%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +4.60%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.03"
$ws.Range("E10").Value = "  -2.87%  "
# Row 11
$ws.Range("E11").Value = "  +1.36%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0830"
$ws.Range("E12").Value = "  -2.19%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.80"
$ws.Range("E13").Value = "  +1.54%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.74"
$ws.Range("E14").Value = "  +2.28%  "
# Row 15
$ws.Range("D15").Value = "3.207.25"
$ws.Range("E15").Value = "  -0.53%  "
# Row 16
$ws.Range("D16").Value = "2.784.12"
$ws.Range("E16").Value = "  -0.63%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.919"
$ws.Range("E17").Value = "  -2.20%  "
# Row 18
$ws.Range("D18").Value = "51.322.49"
$ws.Range("E18").Value = "  -1.14%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +2.76%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.10"
$ws.Range("E20").Value = "  -0.09%  "
# Row 21
$ws.Range("E21").Value = "  +1.21%  "
# Row 22
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  -1.71%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.13"
$ws.Range("E23").Value = "  -0.40%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.36"
$ws.Range("E24").Value = "  -1.76%  "
# Row 25
$ws.Range("E25").Value = "  -0.56%  "
# Row 26
$ws.Range("E26").Value = "  -0.07%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.76"
$ws.Range("E27").Value = "  -3.04%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  +0.89%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.20"
$ws.Range("E29").Value = "  -0.79%  "
# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  -2.74%  "
# Row 31
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.75"
$ws.Range("E31").Value = "  +7.50%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  +8.64%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.85"
$ws.Range("E33").Value = "  -0.19%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0442"
$ws.Range("E34").Value = "  -6.04%  "
# Row 35
$ws.Range("E35").Value = "  +5.02%  "
# Row 36
$ws.Range("E36").Value = "  -0.04%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0837"
$ws.Range("E37").Value = "  -1.08%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.26"
$ws.Range("E38").Value = "  -3.66%  "
# Row 39
$ws.Range("E39").Value = "  -4.24%  "
# Row 40
$ws.Range("E40").Value = "  -2.43%  "
# Row 41
$ws.Range("E41").Value = "  -1.45%  "
# Row 42
$ws.Range("E42").Value = "  -4.31%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.07"
$ws.Range("E43").Value = "  +0.26%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.18"
$ws.Range("E44").Value = "  -2.62%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.70"
$ws.Range("E45").Value = "  -1.46%  "
# Row 46
$ws.Range("D46").Value = "2.125.10"
$ws.Range("E46").Value = "  +1.87%  "
# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("E47").Value = "  +5.79%  "
# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.27"
$ws.Range("E48").Value = "  +0.29%  "
# Row 49
$ws.Range("E49").Value = "  +18.07%  "
# Row 50
$ws.Range("E50").Value = "  -6.15%  "
# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.34"
$ws.Range("E51").Value = "  +8.63%  "
